$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.874.22"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.434.96"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.43"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.32"
$ws.Range("E6").Value = "  -0.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.433.58"
$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.477"
$ws.Range("E9").Value = "  -0.11%  "

$ws.Range("E10").Value = "  +1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("E11").Value = "  -1.53%  "

$ws.Range("E12").Value = "  +2.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.022.16"
$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("E14").Value = "  +2.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.94"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.449.76"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("E17").Value = "  -1.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.877.05"
$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("E19").Value = "  +1.94%  "

$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.16"
$ws.Range("E21").Value = "  -1.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.33"
$ws.Range("E22").Value = "  -3.22%  "

$ws.Range("E23").Value = "  -0.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.38"
$ws.Range("E24").Value = "  -1.48%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.597.49"
$ws.Range("E26").Value = "  +0.65%  "

$ws.Range("E27").Value = "  -3.88%  "

$ws.Range("E28").Value = "  -5.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.56"
$ws.Range("E29").Value = "  -1.50%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.08"
$ws.Range("E31").Value = "  -1.30%  "

$ws.Range("E32").Value = "  -2.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("E34").Value = "  -2.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.31"
$ws.Range("E35").Value = "  -7.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.27"
$ws.Range("E36").Value = "  -1.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.06"
$ws.Range("E37").Value = "  -0.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.58"
$ws.Range("E38").Value = "  -1.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.50"
$ws.Range("E39").Value = "  +6.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "168.65"
$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.469.57"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0764"
$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.30"
$ws.Range("E44").Value = "  -1.54%  "

$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("E46").Value = "  -0.37%  "

$ws.Range("E47").Value = "  -3.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.576.27"
$ws.Range("E48").Value = "  +1.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.26"
$ws.Range("E49").Value = "  +3.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.78"
$ws.Range("E50").Value = "  +0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.65"
$ws.Range("E51").Value = "  -4.09%  "
